$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "63.015.31"
Set-TextCell "E2" "  -0.34%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.606.35"
Set-TextCell "E3" "  -1.53%  "

# Row 4 - TetherUSD
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.16%  "

# Row 5 - BNB
Set-TextCell "D5" "603.41"
Set-TextCell "E5" "  +1.75%  "

# Row 6 - Solana
Set-TextCell "D6" "145.31"
Set-TextCell "E6" "  +1.48%  "

# Row 7 - USDC
Set-TextCell "D7" "0.999"
Set-TextCell "E7" "  -0.16%  "

# Row 8 - XRP
Set-TextCell "D8" "0.584"
Set-TextCell "E8" "  -0.38%  "

# Row 9 - LidoStakedEther
Set-TextCell "D9" "2.604.65"
Set-TextCell "E9" "  -1.47%  "

# Row 10 - Dogecoin
Set-TextCell "D10" "0.108"
Set-TextCell "E10" "  +1.34%  "

# Row 11 - Toncoin
Set-TextCell "D11" "5.49"
Set-TextCell "E11" "  -3.07%  "

# Row 12 - Cardano
Set-TextCell "D12" "0.369"
Set-TextCell "E12" "  +4.00%  "

# Row 13 - TRON
Set-TextCell "E13" "  -0.26%  "

# Row 14 - Avalanche
Set-TextCell "D14" "27.11"
Set-TextCell "E14" "  -0.94%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextCell "D15" "3.070.78"
Set-TextCell "E15" "  -1.79%  "

# Row 16 - WrappedBTC
Set-TextCell "D16" "62.766.86"
Set-TextCell "E16" "  -0.70%  "

# Row 17 - ShibaInu
Set-TextCell "E17" "  +1.12%  "

# Row 18 - WrappedEther
Set-TextCell "D18" "2.602.39"
Set-TextCell "E18" "  -1.48%  "

# Row 19 - Chainlink
Set-TextCell "D19" "11.42"
Set-TextCell "E19" "  +0.19%  "

# Row 20 - Polkadot
Set-TextCell "D20" "4.51"
Set-TextCell "E20" "  +3.07%  "

# Row 21 - BitcoinCash
Set-TextCell "D21" "341.40"
Set-TextCell "E21" "  +0.62%  "

# Row 22 - Uniswap
Set-TextCell "D22" "6.83"
Set-TextCell "E22" "  +0.77%  "

# Row 23 - Dai
Set-TextCell "E23" "  -0.05%  "

# Row 24 - LEO
Set-TextCell "D24" "5.70"
Set-TextCell "E24" "  -1.52%  "

# Row 25 - Litecoin
Set-TextCell "D25" "65.98"
Set-TextCell "E25" "  -1.94%  "

# Row 26 - Fetch.AI
Set-TextCell "D26" "1.68"
Set-TextCell "E26" "  +0.63%  "

# Row 27 - was InternetComputer(DFINITY), now SuiNetwork
Set-TextCell "B27" "SuiNetwork"
Set-TextCell "C27" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextCell "D27" "1.58"
Set-TextCell "E27" "  +3.71%  "

# Row 28 - was SuiNetwork, now InternetComputer(DFINITY)
Set-TextCell "B28" "InternetComputer(DFINITY)"
Set-TextCell "C28" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D28" "8.96"
Set-TextCell "E28" "  +6.23%  "

# Row 29 - Bittensor
Set-TextCell "D29" "553.30"
Set-TextCell "E29" "  +3.72%  "

# Row 30 - Binance-PegBSC-USD
Set-TextCell "E30" "  -0.05%  "

# Row 31 - Kaspa
Set-TextCell "E31" "  -2.98%  "

# Row 32 - Aptos
Set-TextCell "E32" "  -1.15%  "

# Row 33 - PancakeSwap
Set-TextCell "D33" "2.02"
Set-TextCell "E33" "  +1.65%  "

# Row 34 - PEPE
Set-TextCell "D34" "0.0₃0841"
Set-TextCell "E34" "  +4.30%  "

# Row 35 - ImmutableX
Set-TextCell "D35" "1.75"
Set-TextCell "E35" "  -5.07%  "

# Row 36 - NEARProtocol
Set-TextCell "D36" "5.12"
Set-TextCell "E36" "  +0.84%  "

# Row 37 - Monero
Set-TextCell "D37" "167.44"
Set-TextCell "E37" "  -3.74%  "

# Row 38 - FirstDigitalUSD
Set-TextCell "D38" "0.998"
Set-TextCell "E38" "  -0.23%  "

# Row 39 - PolygonEcosystemToken
Set-TextCell "D39" "0.401"
Set-TextCell "E39" "  -1.26%  "

# Row 40 - Stacks
Set-TextCell "E40" "  +5.37%  "

# Row 41 - EthereumClassic
Set-TextCell "D41" "18.91"
Set-TextCell "E41" "  -0.64%  "

# Row 42 - USDe
Set-TextCell "E42" "  -0.06%  "

# Row 43 - Aave
Set-TextCell "D43" "164.60"
Set-TextCell "E43" "  -4.53%  "

# Row 44 - OKB
Set-TextCell "D44" "39.53"
Set-TextCell "E44" "  -1.37%  "

# Row 45 - Filecoin
Set-TextCell "D45" "3.73"
Set-TextCell "E45" "  -0.60%  "

# Row 46 - was Hedera, now InjectiveProtocol
Set-TextCell "B46" "InjectiveProtocol"
Set-TextCell "C46" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D46" "21.63"
Set-TextCell "E46" "  -1.96%  "

# Row 47 - was InjectiveProtocol, now Hedera
Set-TextCell "B47" "Hedera"
Set-TextCell "C47" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D47" "0.0561"
Set-TextCell "E47" "  +0.38%  "

# Row 48 - Mantle
Set-TextCell "E48" "  -1.52%  "

# Row 49 - VeChain
Set-TextCell "D49" "0.0244"
Set-TextCell "E49" "  +1.70%  "

# Row 50 - Stellar
Set-TextCell "D50" "0.0953"
Set-TextCell "E50" "  -0.65%  "

# Row 51 - dogwifhat
Set-TextCell "E51" "  +11.68%  "
